# Opcion para restablecer concurso
# Add a new "RONDA EMPATE" header column to the "Puntuaciones Detalle" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Puntuaciones Detalle")

# Copy the style of the last existing header cell (H1) onto the new one (I1)
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats

# New header text for the added column
$ws.Range("I1").Value = "RONDA EMPATE"
